# wrap <fn>if, countif, averageif, concatenate
$wb = $excel.ActiveWorkbook

# --- Sheet1: selection moves from B4 to B1 ---
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws1.Range("B1").Select()

# --- Sheet3: rework SUMIF in B4, add AVERAGE/AVERAGEIF/COUNTIF/CONCATENATE rows ---
$ws3 = $wb.Worksheets.Item("Sheet3")

# Drop the sum_range argument so SUMIF returns the criteria match itself
$ws3.Range("B4").Formula = "=SUMIF(A1:A3,"">""&A1)"

$ws3.Range("B10").Formula = "=AVERAGE(A1:A3)"
$ws3.Range("B11").Formula = "=AVERAGEIF(A1:A3,"">100"")"
$ws3.Range("B12").Formula = "=COUNTIF(A1:A3,"">100"")"
$ws3.Range("B13").Formula = "=CONCATENATE(B1,"":"",B2,"":"",B3)"

# Selection ends up on the newly added last cell
$ws3.Range("B13").Select()
